$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (index 1) - update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 151
$ws1.Range("F4").Value = 460
$ws1.Range("F5").Value = 1794
$ws1.Range("F7").Value = 2261
$ws1.Range("F11").Value = 5116
$ws1.Range("F17").Value = 212
$ws1.Range("F20").Value = 128
$ws1.Range("F21").Value = 4160
$ws1.Range("F22").Value = 746
$ws1.Range("F23").Value = 757
$ws1.Range("F24").Value = 37
$ws1.Range("F27").Value = 136
$ws1.Range("F33").Value = 31
$ws1.Range("F34").Value = 1053
$ws1.Range("F35").Value = 7
$ws1.Range("F36").Value = 2672
$ws1.Range("F38").Value = 63

# Sheet 4: "全部类型" (index 4) - same updates, row numbers shifted by 1
# from row 34 onward because of an extra row (row 33) present only here
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 151
$ws4.Range("F4").Value = 460
$ws4.Range("F5").Value = 1794
$ws4.Range("F7").Value = 2261
$ws4.Range("F11").Value = 5116
$ws4.Range("F17").Value = 212
$ws4.Range("F20").Value = 128
$ws4.Range("F21").Value = 4160
$ws4.Range("F22").Value = 746
$ws4.Range("F23").Value = 757
$ws4.Range("F24").Value = 37
$ws4.Range("F27").Value = 136
$ws4.Range("F34").Value = 31
$ws4.Range("F35").Value = 1053
$ws4.Range("F36").Value = 7
$ws4.Range("F37").Value = 2672
$ws4.Range("F39").Value = 63
